$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "merging co-cooking with the spaces" -- remove the three separate
# co_cooking_A / co_cooking_B / co_cooking_C rows (original rows 3, 5, 7).
# After deleting row 3 (co_cooking_A), co_cooking_B (was row 5) shifts to
# row 4; after deleting that, co_cooking_C (was row 7) shifts to row 5.
$ws.Rows("3").Delete()
$ws.Rows("4").Delete()
$ws.Rows("5").Delete()

# Update the remaining voxel_total figures for the first three spaces and
# add a new voxel_depth column (only populated for student_housing).
$ws.Range("D1").Value = "voxel_depth"
$ws.Range("B2").Value = 1879
$ws.Range("D2").Value = 13
$ws.Range("B3").Value = 1948
$ws.Range("B4").Value = 2355

# New columns B:D now hold real content -- give them explicit widths.
$ws.Columns("B").ColumnWidth = 13.5
$ws.Columns("C").ColumnWidth = 15.35
$ws.Columns("D").ColumnWidth = 16.45

# Leave the selection where the editor last clicked.
$ws.Range("M15").Select() | Out-Null
